# Actualización automática 2025-08-14 17:20:08
# Updates the August sales figure for client "TAMAYO CONDO LUIS ALFREDO"
# (advisor ALMEIDA CUATIN JHONATHANN CARLOS) in the PORCELANATO group from
# 0 to 1081.02, and propagates the change through the dependent
# summary/rollup cells on the other two sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": per-client totals by product group ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
# M29 = PORCELANATO sales for TAMAYO CONDO LUIS ALFREDO
$wsGrupo.Range("M29").Value = 1081.02
# M34 = count of advisors meeting the PORCELANATO group target
$wsGrupo.Range("M34").Value = "3 de 32"

# --- Sheet "VENTA MENSUAL": per-client totals by month ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
# F29 = agosto (August) sales for TAMAYO CONDO LUIS ALFREDO
$wsMensual.Range("F29").Value = 1081.02
# F34 = TOTAL agosto (August) sales across all clients
$wsMensual.Range("F34").Value = 9513.66

# --- Sheet "CUMPLIMIENTO MENSUAL": compliance rollups ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Row 16 = PORCELANATO group totals
$wsCumpl.Range("D16").Value = 4819.17
$wsCumpl.Range("E16").Value = 17053.93
$wsCumpl.Range("F16").Value = 0.2203240510032872
# Row 19 = grand TOTAL across all groups
$wsCumpl.Range("D19").Value = 9633.74
$wsCumpl.Range("E19").Value = 22475.54107555787
$wsCumpl.Range("F19").Value = 0.3000297632740636
